# Add a new paragraph ("This is the next Change") right after the
# existing last paragraph of the document (mirrors the author's edit
# of appending a new paragraph following the "...VNRVJIETs" paragraph).

$d = $word.ActiveDocument

# Locate the last paragraph in the document body and insert a brand new
# paragraph mark right after it; this yields a new empty paragraph that
# inherits the run/paragraph formatting (en-US language) from its
# neighbour, matching how Word itself creates a new paragraph.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

# The newly created paragraph is now the last paragraph in the document;
# give it the requested text.
$d.Paragraphs.Last.Range.Text = "This is the next Change"
